$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 45

$ws.Range("E17").Value = 108

$ws.Range("E18").Value = 104

$ws.Range("E23").Value = 4

$ws.Range("F28").Value = 11
$ws.Range("H28").Value = 11

$ws.Range("F36").Value = 41
$ws.Range("H36").Value = 41

$ws.Range("E37").Value = 47
$ws.Range("F37").Value = 26
$ws.Range("H37").Value = 26

$ws.Range("E42").Value = 34

$ws.Range("F44").Value = 14
$ws.Range("H44").Value = 14

$ws.Range("F50").Value = 8
$ws.Range("H50").Value = 8

$ws.Range("F60").Value = 9
$ws.Range("H60").Value = 9

$ws.Range("E62").Value = 39

$ws.Range("F65").Value = 9
$ws.Range("H65").Value = 9

$ws.Range("E71").Value = 31

$ws.Range("E76").Value = 49

$ws.Range("F77").Value = 19
$ws.Range("H77").Value = 19

$ws.Range("E78").Value = 44
$ws.Range("F78").Value = 18
$ws.Range("H78").Value = 18

$ws.Range("E80").Value = 23
$ws.Range("F80").Value = 8
$ws.Range("H80").Value = 8

$wb.Save()
